$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# New rows 10 & 11: clone formatting+values from row 4 (same shape) first,
# while B4 still reads "Yes" - then adjust the TestName (A) and Browser
# Version (D) cells for the two new rows.
$ws.Range("A4:I4").Copy($ws.Range("A10:I10"))
$ws.Range("A4:I4").Copy($ws.Range("A11:I11"))

$ws.Range("A10").Value = "verifyAmazonHamburgerSubMenuPageTitle1"
$ws.Range("D10").Value = "'84.0"
$ws.Range("E10").Value = "'"
$ws.Range("F10").Value = "'"

$ws.Range("A11").Value = "verifyAmazonHamburgerSubMenuPageTitle2"
$ws.Range("D11").Value = "'84.0"
$ws.Range("E11").Value = "'"
$ws.Range("F11").Value = "'"

$ws.Rows(10).RowHeight = 18
$ws.Rows(11).RowHeight = 18

# B4:B9 flip from "Yes" to "No" (Execute column)
$ws.Range("B4:B9").Value = "No"

$ws.Range("B12").Select()
